$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.327.04"
$ws.Range("E2").Value = "  +2.54%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.894.59"
$ws.Range("E3").Value = "  +0.72%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  -0.24%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.63"
$ws.Range("E5").Value = "  -0.14%  "

$ws.Range("E6").Value = "  -0.36%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5140"
$ws.Range("E7").Value = "  +0.50%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3929"
$ws.Range("E8").Value = "  -0.83%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08426"
$ws.Range("E9").Value = "  -0.10%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.37"
$ws.Range("E10").Value = "  +1.66%  "

$ws.Range("E11").Value = "  +0.21%  "

$ws.Range("E12").Value = "  -0.07%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.893.49"
$ws.Range("E13").Value = "  +0.77%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.72"
$ws.Range("E14").Value = "  +0.95%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.295"
$ws.Range("E15").Value = "  +0.28%  "

$ws.Range("E16").Value = "  -0.31%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.23"
$ws.Range("E17").Value = "  +2.25%  "

$ws.Range("E18").Value = "  -0.21%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06739"
$ws.Range("E19").Value = "  +0.00%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.86"
$ws.Range("E20").Value = "  +0.75%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.003"
$ws.Range("E21").Value = "  -0.40%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.012"
$ws.Range("E22").Value = "  +0.75%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "29.342.04"
$ws.Range("E23").Value = "  +2.53%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.16"
$ws.Range("E24").Value = "  +0.31%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.213"
$ws.Range("E25").Value = "  -1.93%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.112.72"
$ws.Range("E26").Value = "  +0.89%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.30"
$ws.Range("E27").Value = "  -1.11%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.86"
$ws.Range("E28").Value = "  +0.11%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.430"
$ws.Range("E29").Value = "  +1.87%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.05"
$ws.Range("E30").Value = "  -0.24%  "

$ws.Range("E31").Value = "  +1.03%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1049"
$ws.Range("E32").Value = "  -0.31%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.122"
$ws.Range("E33").Value = "  +5.62%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.656"
$ws.Range("E34").Value = "  +1.34%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02483"
$ws.Range("E35").Value = "  +1.08%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06582"
$ws.Range("E36").Value = "  +0.84%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2202"
$ws.Range("E37").Value = "  +0.33%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "9.042"
$ws.Range("E38").Value = "  +1.22%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.202"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.230"
$ws.Range("E40").Value = "  +2.63%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6519"
$ws.Range("E41").Value = "  +1.13%  "

$ws.Range("E42").Value = "  -2.28%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.26"
$ws.Range("E43").Value = "  +0.69%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6068"
$ws.Range("E44").Value = "  -0.08%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.17"
$ws.Range("E45").Value = "  +1.05%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.673"
$ws.Range("E46").Value = "  -0.61%  "

$ws.Range("E47").Value = "  +1.82%  "

$ws.Range("E48").Value = "  +1.46%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "123.49"
$ws.Range("E49").Value = "  +0.68%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.155"
$ws.Range("E50").Value = "  -3.40%  "

$ws.Range("E51").Value = "  +0.70%  "
